$wb = $excel.ActiveWorkbook

# Updated "想去人数" (interest count) figures for the 展览 (Exhibitions)
# and 全部类型 (All types) sheets, as of the data refresh at 456a3b4.
$updates = @{
    "F5"  = 2684
    "F9"  = 1420
    "F13" = 1209
    "F15" = 364
    "F19" = 107
    "F21" = 89
    "F22" = 2605
    "F23" = 40
    "F24" = 298
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
